# Updates the result-line active power values (pl_mw) for the 380 kV case (Case_2_3).
# Each entry is (row, column index, new Value2) matching cells B:N (skipping D/K which stay 0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 2, 1.235034407920523),
    @(2, 3, 0.1905543085327963),
    @(2, 5, 0.1750444730758582),
    @(2, 6, 2.083185004842619),
    @(2, 7, 1.047729313344064),
    @(2, 8, 1.036581200750263),
    @(2, 9, 0.9968959980320093),
    @(2, 10, 0.0519340307173195),
    @(2, 12, 0.442840988115222),
    @(2, 13, 0.3493261213143342),
    @(2, 14, 1.589797112282355),
    @(3, 2, 1.148705948682164),
    @(3, 3, 0.1728947131852294),
    @(3, 5, 0.1754753712506512),
    @(3, 6, 2.080451053136699),
    @(3, 7, 1.041538780712003),
    @(3, 8, 1.039183797146976),
    @(3, 9, 1.002348751589373),
    @(3, 10, 0.05162651519041006),
    @(3, 12, 0.4390765830141348),
    @(3, 13, 0.3349726243342204),
    @(3, 14, 1.608120184331),
    @(4, 2, 1.09611105637515),
    @(4, 3, 0.161969185187786),
    @(4, 5, 0.1757613576528221),
    @(4, 6, 2.079867179875876),
    @(4, 7, 1.0384837223868),
    @(4, 8, 1.041325294943988),
    @(4, 9, 1.006254125864352),
    @(4, 10, 0.05144132818572267),
    @(4, 12, 0.4369183537759227),
    @(4, 13, 0.3262968858204474),
    @(4, 14, 1.619955938006633),
    @(5, 2, 1.074782480408885),
    @(5, 3, 0.1574961518726923),
    @(5, 5, 0.1758832976927072),
    @(5, 6, 2.079904553116435),
    @(5, 7, 1.037425931326624),
    @(5, 8, 1.042334573854234),
    @(5, 9, 1.007985617787284),
    @(5, 10, 0.051366790386421),
    @(5, 12, 0.4360775079904826),
    @(5, 13, 0.3227961960998016),
    @(5, 14, 1.624926255026706),
    @(6, 2, 1.071247210539696),
    @(6, 3, 0.156752150729119),
    @(6, 5, 0.1759038721906889),
    @(6, 6, 2.079927387225581),
    @(6, 7, 1.037261578693546),
    @(6, 8, 1.042510412371499),
    @(6, 9, 1.008281583714201),
    @(6, 10, 0.05135446989436687),
    @(6, 12, 0.435940224814189),
    @(6, 13, 0.3222170136789586),
    @(6, 14, 1.625760459141528),
    @(7, 2, 1.095822988378274),
    @(7, 3, 0.1619089445533746),
    @(7, 5, 0.1757629803026965),
    @(7, 6, 2.079866569184063),
    @(7, 7, 1.038468699342701),
    @(7, 8, 1.041338353428785),
    @(7, 9, 1.006276910612669),
    @(7, 10, 0.05144031916830549),
    @(7, 12, 0.4369068571504755),
    @(7, 13, 0.3262495333482818),
    @(7, 14, 1.620022373672711),
    @(8, 2, 1.205183527392762),
    @(8, 3, 0.1844823573582346),
    @(8, 5, 0.1751886106824898),
    @(8, 6, 2.082015117169689),
    @(8, 7, 1.045439745168736),
    @(8, 8, 1.037365732739616),
    @(8, 9, 0.9986603457161678),
    @(8, 10, 0.05182725501745011),
    @(8, 12, 0.4415113233286689),
    @(8, 13, 0.3443486410047996),
    @(8, 14, 1.595993408460522),
    @(9, 2, 1.422874175558434),
    @(9, 3, 0.2280989175822299),
    @(9, 5, 0.1742315750398085),
    @(9, 6, 2.094917559057336),
    @(9, 7, 1.065049439773745),
    @(9, 8, 1.033891771552078),
    @(9, 9, 0.988153347736251),
    @(9, 10, 0.05261422037499131),
    @(9, 12, 0.4517501880494592),
    @(9, 13, 0.380924119746112),
    @(9, 14, 1.553517640347646),
    @(10, 2, 1.584764556443474),
    @(10, 3, 0.2597569276090042),
    @(10, 5, 0.1736308431901961),
    @(10, 6, 2.109702176566245),
    @(10, 7, 1.083109607757621),
    @(10, 8, 1.033977809567943),
    @(10, 9, 0.9831440833646781),
    @(10, 10, 0.05320884087428723),
    @(10, 12, 0.4600038686192107),
    @(10, 13, 0.4084505147297648),
    @(10, 14, 1.525142555705948),
    @(11, 2, 1.658833887082039),
    @(11, 3, 0.274077049468616),
    @(11, 5, 0.1733796215289871),
    @(11, 6, 2.117582122979385),
    @(11, 7, 1.092125977992964),
    @(11, 8, 1.034591485133774),
    @(11, 9, 0.9814560168879325),
    @(11, 10, 0.05348276601887747),
    @(11, 12, 0.4639162292526606),
    @(11, 13, 0.4211140276191117),
    @(11, 14, 1.512848995738988),
    @(12, 2, 1.686942441894757),
    @(12, 3, 0.279488126603411),
    @(12, 5, 0.1732876484869186),
    @(12, 6, 2.120732135352142),
    @(12, 7, 1.095655897321421),
    @(12, 8, 1.034906597888551),
    @(12, 9, 0.9809019035974487),
    @(12, 10, 0.05358697352059849),
    @(12, 12, 0.465420291977523),
    @(12, 13, 0.4259295788712336),
    @(12, 14, 1.50828214903548),
    @(13, 2, 1.6808861034595),
    @(13, 3, 0.2783232714170367),
    @(13, 5, 0.1733073162164844),
    @(13, 6, 2.120046337940806),
    @(13, 7, 1.094890516427682),
    @(13, 8, 1.034835051645899),
    @(13, 9, 0.9810174531116687),
    @(13, 10, 0.05356450955127201),
    @(13, 12, 0.4650953651097609),
    @(13, 13, 0.4248915707484358),
    @(13, 14, 1.509261765258483),
    @(14, 2, 1.661145195584083),
    @(14, 3, 0.2745224551165961),
    @(14, 5, 0.1733719916068504),
    @(14, 6, 2.117837948238034),
    @(14, 7, 1.092414066748518),
    @(14, 8, 1.034615751042026),
    @(14, 9, 0.9814087226542654),
    @(14, 10, 0.05349132973928405),
    @(14, 12, 0.4640395185480628),
    @(14, 13, 0.4215098032108884),
    @(14, 14, 1.512471505062113),
    @(15, 2, 1.649061119141322),
    @(15, 3, 0.2721928305302299),
    @(15, 5, 0.1734120182050063),
    @(15, 6, 2.116506871704885),
    @(15, 7, 1.090912240486205),
    @(15, 8, 1.0344921998811),
    @(15, 9, 0.9816594772393827),
    @(15, 10, 0.0534465667933901),
    @(15, 12, 0.4633957127386878),
    @(15, 13, 0.4194409913250396),
    @(15, 14, 1.514449086410991),
    @(16, 2, 1.579932414065752),
    @(16, 3, 0.2588194424176322),
    @(16, 5, 0.173647703856278),
    @(16, 6, 2.109210436551223),
    @(16, 7, 1.082536519613399),
    @(16, 8, 1.033949274892791),
    @(16, 9, 0.9832663039270031),
    @(16, 10, 0.05319100693177958),
    @(16, 12, 0.4597513481191129),
    @(16, 13, 0.4076257560573211),
    @(16, 14, 1.525958341218921),
    @(17, 2, 1.537632262738157),
    @(17, 3, 0.2505945188129886),
    @(17, 5, 0.1737979292500018),
    @(17, 6, 2.105030004245805),
    @(17, 7, 1.077603704247849),
    @(17, 8, 1.033763437984447),
    @(17, 9, 0.9844034552583736),
    @(17, 10, 0.05303509675379559),
    @(17, 12, 0.4575559422010542),
    @(17, 13, 0.4004136133824687),
    @(17, 14, 1.533176335198739),
    @(18, 2, 1.513342354712961),
    @(18, 3, 0.245856108553852),
    @(18, 5, 0.1738864114990966),
    @(18, 6, 2.102734179702026),
    @(18, 7, 1.074841817395736),
    @(18, 8, 1.033710620938535),
    @(18, 9, 0.9851131005498601),
    @(18, 10, 0.05294574502138261),
    @(18, 12, 0.4563080566625786),
    @(18, 13, 0.396278721934749),
    @(18, 14, 1.537385758755359),
    @(19, 2, 1.505125111696429),
    @(19, 3, 0.2442504487372332),
    @(19, 5, 0.1739167270966426),
    @(19, 6, 2.101975512022349),
    @(19, 7, 1.073919614916747),
    @(19, 8, 1.033702021509711),
    @(19, 9, 0.9853629156435062),
    @(19, 10, 0.05291554811545751),
    @(19, 12, 0.4558880996027455),
    @(19, 13, 0.3948810167215342),
    @(19, 14, 1.538820925470091),
    @(20, 2, 1.542131050997455),
    @(20, 3, 0.2514708675856525),
    @(20, 5, 0.1737817226720892),
    @(20, 6, 2.105463773192966),
    @(20, 7, 1.078121010219803),
    @(20, 8, 1.033777623657826),
    @(20, 9, 0.9842766492994883),
    @(20, 10, 0.05305166026877117),
    @(20, 12, 0.4577881106286554),
    @(20, 13, 0.401179978937428),
    @(20, 14, 1.532401982055212),
    @(21, 2, 1.666941955000539),
    @(21, 3, 0.275639162612805),
    @(21, 5, 0.1733529092353514),
    @(21, 6, 2.118482098904664),
    @(21, 7, 1.093138318485671),
    @(21, 8, 1.034677918873854),
    @(21, 9, 0.9812914856866897),
    @(21, 10, 0.05351281157149401),
    @(21, 12, 0.4643490357782696),
    @(21, 13, 0.4225025650191085),
    @(21, 14, 1.51152632508823),
    @(22, 2, 1.748863008160242),
    @(22, 3, 0.2913667013466466),
    @(22, 5, 0.1730910632306415),
    @(22, 6, 2.127958206259805),
    @(22, 7, 1.103627115000563),
    @(22, 8, 1.035748562847061),
    @(22, 9, 0.9798367438579305),
    @(22, 10, 0.05381697978043931),
    @(22, 12, 0.4687682744757637),
    @(22, 13, 0.4365554820781767),
    @(22, 14, 1.498398528223852),
    @(23, 2, 1.705108495183936),
    @(23, 3, 0.2829788106625983),
    @(23, 5, 0.1732291350067401),
    @(23, 6, 2.12281205033743),
    @(23, 7, 1.097967207944833),
    @(23, 8, 1.035132977611482),
    @(23, 9, 0.9805677006884181),
    @(23, 10, 0.05365439015913154),
    @(23, 12, 0.4663976773137648),
    @(23, 13, 0.4290445053584548),
    @(23, 14, 1.505357861873019),
    @(24, 2, 1.540097058381662),
    @(24, 3, 0.251074700815991),
    @(24, 5, 0.173789043075872),
    @(24, 6, 2.105267331202754),
    @(24, 7, 1.077886905547317),
    @(24, 8, 1.033771042036435),
    @(24, 9, 0.984333804212902),
    @(24, 10, 0.05304417101990921),
    @(24, 12, 0.4576831027843156),
    @(24, 13, 0.4008334693403341),
    @(24, 14, 1.532751881139859),
    @(25, 2, 1.363638752897828),
    @(25, 3, 0.2163681147729051),
    @(25, 5, 0.1744724377555089),
    @(25, 6, 2.090496047438236),
    @(25, 7, 1.059105326749631),
    @(25, 8, 1.034368758893635),
    @(25, 9, 0.990520563068948),
    @(25, 10, 0.0523983854660024),
    @(25, 12, 0.4488514563741433),
    @(25, 13, 0.3709141400980229),
    @(25, 14, 1.564511187666621)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value2 = $u[2]
}

